# Insert a new weekly data row at row 387 (shifting all subsequent rows down by one),
# and populate it with the new "Paine" / "1a (guarda)" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 387; everything currently at
# row 387 and below (up to 458) moves down to 388..459.
$ws.Rows.Item(387).Insert()

# Populate the newly inserted row 387 with the new record's values.
$ws.Cells.Item(387, 1).Value  = 11
$ws.Cells.Item(387, 2).Value  = 'Vega Monumental Concepción'
$ws.Cells.Item(387, 3).Value  = 'Bíobío'
$ws.Cells.Item(387, 4).Value  = 45211
$ws.Cells.Item(387, 5).Value  = 8
$ws.Cells.Item(387, 6).Value  = 100112045
$ws.Cells.Item(387, 7).Value  = 'Zapallo'
$ws.Cells.Item(387, 8).Value  = 'Paine'
$ws.Cells.Item(387, 9).Value  = '1a (guarda)'
$ws.Cells.Item(387, 10).Value = 500
$ws.Cells.Item(387, 11).Value = 350
$ws.Cells.Item(387, 12).Value = 350
$ws.Cells.Item(387, 13).Value = 350
$ws.Cells.Item(387, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(387, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(387, 16).Value = 350
$ws.Cells.Item(387, 17).Value = 1
$ws.Cells.Item(387, 18).Value = 'Hortaliza'
